$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.435.65"
$ws.Range("E2").Value = "  +0.08%  "

$ws.Range("D3").Value = "3.689.77"
$ws.Range("E3").Value = "  -0.04%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "679.58"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.98%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "161.17"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").Value = "  +0.01%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.147"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.00%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "7.18"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -0.76%  "

$ws.Range("E11").Value = "  +0.33%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.0000234"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.16%  "

$ws.Range("D13").Value = "4.311.50"
$ws.Range("E13").Value = "  -0.08%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "32.47"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -0.25%  "

$ws.Range("D15").Value = "3.692.13"
$ws.Range("E15").Value = "  +0.17%  "

$ws.Range("D16").Value = "69.401.09"
$ws.Range("E16").Value = "  -0.05%  "

$ws.Range("E17").Value = "  +2.68%  "

$ws.Range("E18").Value = "  +0.55%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "6.48"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +0.39%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "471.65"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.64%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "9.81"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -1.04%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.651"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.48%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "80.31"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.87%  "

$ws.Range("D24").Value = "3.836.01"
$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("E25").Value = "  -0.08%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.0000126"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.26%  "

$ws.Range("E27").Value = "  -1.55%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "9.15"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -0.69%  "

$ws.Range("E29").Value = "  -0.24%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.75"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -0.97%  "

$ws.Range("E31").Value = "  -0.56%  "

$ws.Range("E32").Value = "  -1.47%  "

$ws.Range("E33").Value = "  +0.17%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "26.97"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +1.00%  "

$ws.Range("D35").Value = "3.679.72"
$ws.Range("E35").Value = "  +0.36%  "

$ws.Range("E36").Value = "  +2.09%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "8.46"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +2.82%  "

$ws.Range("E38").Value = "  +1.90%  "

$ws.Range("E40").Value = "  -0.31%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("E42").Value = "  -0.23%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "168.61"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +1.62%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.942"
$cell.Style = "Normal"

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "46.69"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -2.51%  "

$ws.Range("E46").Value = "  +0.12%  "

$ws.Range("E47").Value = "  +1.85%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "28.13"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -1.08%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.28"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -2.38%  "

$ws.Range("E50").Value = "  -2.68%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "7.89"
$cell.Style = "Normal"
